# The commit unifies the naming of "Property1" into "DataNode" to match the
# shared DataNode / DataTable / Entity concept used elsewhere in the project.
# The only content-level change is renaming the sole worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"
